$wb = $excel.ActiveWorkbook

# Populate the header row (Name / Salary) on the six previously-empty
# report sheets: TempHelp, Overtime, Retirement, SocialSecurity,
# MedicalAndLifeIns, IndustIns (sheets index 3-8).
for ($i = 3; $i -le 8; $i++) {
    $sheet = $wb.Worksheets.Item($i)
    $sheet.Range("A1").Value = "Name"
    $sheet.Range("B1").Value = "Salary"
    $sheet.Range("A1:B1").Select()
}

# Finalize the report: make Sheet1 the active/selected tab instead of
# Salary (sheet 2).
$wb.Worksheets.Item(1).Activate()
